# Updated symbol list (price + 1h volume %) refresh
# Source: scheduled GitHub Actions data pull (coinranking.com), run timestamp
# Mon Jan 16 10:52:04 UTC 2023
#
# D/E columns store scraped price & percentage strings as literal text
# (e.g. "298.43", "0.68%"); a leading apostrophe forces Excel to keep
# them as text instead of re-interpreting as Number/Percentage,
# matching the existing inline-string cells in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.43"
$ws.Range("E2").Value = "'0.68%"

$ws.Range("D3").Value = "'31.33"
$ws.Range("E3").Value = "'0.62%"

$ws.Range("D4").Value = "'5.140"
$ws.Range("E4").Value = "'0.49%"

$ws.Range("D5").Value = "'0.07946"
$ws.Range("E5").Value = "'7.82%"

$ws.Range("D6").Value = "'2.626"
$ws.Range("E6").Value = "'61.10%"

$ws.Range("D7").Value = "'7.830"
$ws.Range("E7").Value = "'1.70%"

$ws.Range("D8").Value = "'3.825"
$ws.Range("E8").Value = "'2.25%"

$ws.Range("D9").Value = "'0.9078"
$ws.Range("E9").Value = "'-1.04%"

$ws.Range("D10").Value = "'0.1732"
$ws.Range("E10").Value = "'3.25%"

$ws.Range("D11").Value = "'0.07154"
$ws.Range("E11").Value = "'-0.06%"

$ws.Range("D12").Value = "'0.08016"
$ws.Range("E12").Value = "'0.20%"

$ws.Range("D13").Value = "'0.03018"
$ws.Range("E13").Value = "'0.93%"

$ws.Range("D14").Value = "'0.09956"
$ws.Range("E14").Value = "'0.53%"

$ws.Range("D15").Value = "'0.001492"
$ws.Range("E15").Value = "'0.07%"

$ws.Range("D16").Value = "'0.006003"
$ws.Range("E16").Value = "'-2.57%"

$ws.Range("D17").Value = "'3.502"
$ws.Range("E17").Value = "'1.53%"

$ws.Range("E18").Value = "'1.13%"

$ws.Range("D19").Value = "'0.3282"
$ws.Range("E19").Value = "'0.30%"

$ws.Range("E20").Value = "'-1.52%"

$ws.Range("D21").Value = "'4.634"
$ws.Range("E21").Value = "'1.92%"

$ws.Range("E22").Value = "'3.31%"

$ws.Range("D23").Value = "'0.04593"
$ws.Range("E23").Value = "'-0.59%"

$ws.Range("E24").Value = "'3.77%"

$ws.Range("D25").Value = "'0.004454"
$ws.Range("E25").Value = "'0.63%"

$ws.Range("D26").Value = "'0.0001181"
$ws.Range("E26").Value = "'-9.01%"

$ws.Range("D27").Value = "'0.0003428"
$ws.Range("E27").Value = "'82.95%"

$ws.Range("D39").Value = "'0.01857"
$ws.Range("E39").Value = "'10.32%"

$ws.Range("D40").Value = "'0.04513"
$ws.Range("E40").Value = "'2.59%"

$ws.Range("D41").Value = "'0.007007"
$ws.Range("E41").Value = "'-2.90%"

$ws.Range("D42").Value = "'0.1344"
$ws.Range("E42").Value = "'1.33%"

$ws.Range("D43").Value = "'0.002242"
$ws.Range("E43").Value = "'4.93%"

$ws.Range("E44").Value = "'-5.72%"

$ws.Range("D45").Value = "'0.00006441"
$ws.Range("E45").Value = "'7.09%"

$ws.Range("E46").Value = "'-0.06%"

$ws.Range("E47").Value = "'-57.44%"

$ws.Range("D48").Value = "'0.006197"
$ws.Range("E48").Value = "'-39.33%"

$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.06%"

$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.01%"
